# Auto-generated Excel COM-interop script
# Applies the 2026-02-08 21:20 meteocat automatic update diff

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Style-format source cell (untouched, keeps original style index 3)
# used to restore formatting on cells whose new value looks like a pure
# percentage ("NN%"), which Excel would otherwise auto-convert to a
# numeric percent value/format instead of literal text.
$formatSource = $ws.Range("H3")

$ws.Range("E2").Value = "2026-02-08 21:18:24"
$ws.Range("I2").Value = "5.8 mm"
$ws.Range("E3").Value = "2026-02-08 21:18:26"
$ws.Range("I3").Value = "2.3 mm"
$ws.Range("O3").Value = "-4.6 °C"
$ws.Range("E4").Value = "2026-02-08 21:18:29"
$ws.Range("H4").Value = "'70%"
$formatSource.Copy() | Out-Null
$ws.Range("H4").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$ws.Range("J4").Value = "1002.4 hPa"
$ws.Range("O4").Value = "10.2 °C"
$ws.Range("E5").Value = "2026-02-08 21:18:31"
$ws.Range("H5").Value = "'90%"
$formatSource.Copy() | Out-Null
$ws.Range("H5").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$ws.Range("E6").Value = "2026-02-08 21:18:34"
$ws.Range("H6").Value = "'66%"
$formatSource.Copy() | Out-Null
$ws.Range("H6").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$ws.Range("J6").Value = "1002.3 hPa"
$ws.Range("E7").Value = "2026-02-08 21:18:36"
$ws.Range("J7").Value = "1002.6 hPa"
$ws.Range("E8").Value = "2026-02-08 21:18:39"
$ws.Range("J8").Value = "1002.6 hPa"
$ws.Range("O8").Value = "9.4 °C"
$ws.Range("E9").Value = "2026-02-08 21:18:41"
$ws.Range("O9").Value = "10.3 °C"
$ws.Range("E10").Value = "2026-02-08 21:18:44"
$ws.Range("E11").Value = "2026-02-08 21:18:46"
$ws.Range("E12").Value = "2026-02-08 21:18:49"
$ws.Range("E13").Value = "2026-02-08 21:18:51"
$ws.Range("J13").Value = "1003.9 hPa"
$ws.Range("E14").Value = "2026-02-08 21:18:54"
$ws.Range("E15").Value = "2026-02-08 21:18:56"
$ws.Range("H15").Value = "'74%"
$formatSource.Copy() | Out-Null
$ws.Range("H15").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$ws.Range("E16").Value = "2026-02-08 21:18:58"
$ws.Range("I16").Value = "3.8 mm"
$ws.Range("E17").Value = "2026-02-08 21:19:01"
$ws.Range("E18").Value = "2026-02-08 21:19:03"
$ws.Range("J18").Value = "1002.7 hPa"
$ws.Range("O18").Value = "10.2 °C"
$ws.Range("E19").Value = "2026-02-08 21:19:06"
$ws.Range("H19").Value = "'90%"
$formatSource.Copy() | Out-Null
$ws.Range("H19").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$ws.Range("E20").Value = "2026-02-08 21:19:08"
$ws.Range("I20").Value = "9.7 mm"
$ws.Range("E21").Value = "2026-02-08 21:19:11"
$ws.Range("J21").Value = "1003.3 hPa"
$ws.Range("O21").Value = "5.4 °C"
$ws.Range("E22").Value = "2026-02-08 21:19:14"
$ws.Range("E23").Value = "2026-02-08 21:19:16"
$ws.Range("I23").Value = "5.7 mm"
$ws.Range("E24").Value = "2026-02-08 21:19:19"
$ws.Range("J24").Value = "1003.9 hPa"
$ws.Range("E25").Value = "2026-02-08 21:19:21"
$ws.Range("E26").Value = "2026-02-08 21:19:24"
$ws.Range("J26").Value = "1001.7 hPa"
$ws.Range("E27").Value = "2026-02-08 21:19:26"
$ws.Range("E28").Value = "2026-02-08 21:19:29"
$ws.Range("J28").Value = "1002.3 hPa"
$ws.Range("E29").Value = "2026-02-08 21:19:31"
$ws.Range("E30").Value = "2026-02-08 21:19:34"
$ws.Range("J30").Value = "1002.7 hPa"
$ws.Range("E31").Value = "2026-02-08 21:19:36"
$ws.Range("J31").Value = "1001.8 hPa"
$ws.Range("N31").Value = "7.5 °C 20:59 TU"
$ws.Range("O31").Value = "9.6 °C"
$ws.Range("E32").Value = "2026-02-08 21:19:39"
$ws.Range("H32").Value = "'91%"
$formatSource.Copy() | Out-Null
$ws.Range("H32").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$ws.Range("E33").Value = "2026-02-08 21:19:41"
$ws.Range("J33").Value = "1003.5 hPa"
$ws.Range("E34").Value = "2026-02-08 21:19:44"
$ws.Range("E35").Value = "2026-02-08 21:19:46"
$ws.Range("J35").Value = "1004.8 hPa"
$ws.Range("E36").Value = "2026-02-08 21:19:48"
$ws.Range("J36").Value = "1002.7 hPa"
$ws.Range("E37").Value = "2026-02-08 21:19:51"
$ws.Range("H37").Value = "'79%"
$formatSource.Copy() | Out-Null
$ws.Range("H37").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$ws.Range("J37").Value = "1003.5 hPa"
$ws.Range("E38").Value = "2026-02-08 21:19:53"
$ws.Range("H38").Value = "'77%"
$formatSource.Copy() | Out-Null
$ws.Range("H38").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$ws.Range("E39").Value = "2026-02-08 21:19:56"
$ws.Range("E40").Value = "2026-02-08 21:19:58"
$ws.Range("J40").Value = "1004.0 hPa"
$ws.Range("O40").Value = "5.7 °C"
$ws.Range("E41").Value = "2026-02-08 21:20:01"
$ws.Range("H41").Value = "'68%"
$formatSource.Copy() | Out-Null
$ws.Range("H41").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$ws.Range("J41").Value = "1002.7 hPa"
$ws.Range("E42").Value = "2026-02-08 21:20:04"
$ws.Range("E43").Value = "2026-02-08 21:20:06"
$ws.Range("E44").Value = "2026-02-08 21:20:08"
$ws.Range("E45").Value = "2026-02-08 21:20:11"
$ws.Range("J45").Value = "1004.8 hPa"
$ws.Range("E46").Value = "2026-02-08 21:20:13"
$ws.Range("H46").Value = "'72%"
$formatSource.Copy() | Out-Null
$ws.Range("H46").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$ws.Range("J46").Value = "1004.4 hPa"
$ws.Range("K46").Value = "7.9 MJ/m2"

$excel.CutCopyMode = $false

